# Pin Brainstorming deej.xlsx - "Deej script now compiles, Buttons and screen work"
#
# Updates the Pin/Button assignment numbers in column D and a couple of
# related cells (I5, N8) on the active sheet, and leaves the active
# selection on D1 (matching the author's last cursor position when saving).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Decrease the assigned pin numbers in column D by one for rows 2-7
$ws.Range("D2").Value = 9
$ws.Range("D3").Value = 8
$ws.Range("D4").Value = 7
$ws.Range("D5").Value = 6
$ws.Range("D6").Value = 5
$ws.Range("D7").Value = 4

# Button pin reassignment
$ws.Range("I5").Value = 10

# Button value reassignment
$ws.Range("N8").Value = 17

# Leave the selection on D1, as saved in the workbook
$ws.Range("D1").Select()
